$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.166.88"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.656.83"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.37"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5264"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2684"
$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06381"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.64"
$ws.Range("E10").Value = "  -1.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07692"
$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.618"
$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("D13").Value = "1.738.96"
$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("D14").Value = "1.884.59"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5642"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "0.0₅8261"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.78"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "26.146.88"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.692"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").Value = "  +0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "190.91"
$ws.Range("E22").Value = "  -4.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.997"
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.60"
$ws.Range("E25").Value = "  +2.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1204"
$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.288"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.07"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.527"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.279"
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.384"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.803"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9509"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.407"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5784"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.983"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8345"
$ws.Range("E42").Value = "  -2.73%  "

$ws.Range("D43").Value = "1.028.12"
$ws.Range("E43").Value = "  -4.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.57"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").Value = "1.794.85"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.60"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  +3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05358"
$ws.Range("E48").Value = "  +4.18%  "

$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.040"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4342"
$ws.Range("E51").Value = "  -1.53%  "
